# fix: troca planilha modelo
# Swap the "Identificador" column validation rule on the "Dados" sheet from a
# plain text-length check to a custom rule that also forbids blank spaces,
# and propagate the corresponding rule description into the "Instruções"
# sheet / shared strings.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Dados")
$ws2 = $wb.Worksheets.Item("Instruções")

# --- "Dados" sheet -----------------------------------------------------

# Replace the F2:F1048576 "O identificador precisa ter até 25 caracteres."
# textLength validation with a custom formula that also rejects spaces.
$rangeF = $ws1.Range("F2:F1048576")
$rangeF.Validation.Delete()
$rangeF.Validation.Add(7, 1, $null, 'AND(LEN(F2)>=1,LEN(F2)<=25,ISERROR(FIND(" ",F2)))')
$rangeF.Validation.ErrorMessage = "O identificador precisa ter até 25 caracteres e sem espaços em branco."

# Incidental column width adjustment on column G recorded alongside the fix.
$ws1.Columns.Item(7).ColumnWidth = 10.75

# --- "Instruções" sheet --------------------------------------------------

# Update the rule description for "Identificador" to mention the new
# "no blank spaces" restriction (adds a new shared string).
$ws2.Range("B6").Value = "Limite de 25 caracteres e sem espaços em branco"

# Clear the stray cell selection that had been saved with the sheet.
$ws2.Range("A1").Select() | Out-Null
$ws1.Activate() | Out-Null
